# Updated convert_to_dict function and plotting functions
# Apply the data updates to the PV specification sheet (row 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4: P_max 400 -> 800
$ws.Range("A4").Value = 800

# F4: n_eff 0.23 -> 0.28
$ws.Range("F4").Value = 0.28

# G4: no_panels 400 -> 5000
$ws.Range("G4").Value = 5000

# H4: length 1.7 -> 2.4
$ws.Range("H4").Value = 2.4

# I4: width 1 -> 1.3 (also switch from integer format to 2-decimal format,
# matching the style used by the adjacent H column)
$ws.Range("I4").NumberFormat = "#,##0.00"
$ws.Range("I4").Value = 1.3

# Column I (9) width: narrow from the default bestFit width down to ~12
$ws.Columns.Item(9).ColumnWidth = 11.14
